$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @("D2", "289.59"),
    @("E2", "0.94%"),
    @("G2", "5"),
    @("E3", "3.99%"),
    @("G3", "5"),
    @("D4", "5.107"),
    @("E4", "3.74%"),
    @("G4", "5"),
    @("D5", "0.06682"),
    @("E5", "1.96%"),
    @("G5", "5"),
    @("E6", "1.53%"),
    @("G6", "5"),
    @("D7", "1.358"),
    @("E7", "0.40%"),
    @("G7", "5"),
    @("D8", "0.9216"),
    @("E8", "1.25%"),
    @("G8", "5"),
    @("D9", "0.1589"),
    @("E9", "1.29%"),
    @("G9", "5"),
    @("D10", "0.06659"),
    @("E10", "-1.38%"),
    @("G10", "5"),
    @("D11", "0.07700"),
    @("E11", "0.72%"),
    @("G11", "5"),
    @("D12", "0.02939"),
    @("E12", "-1.41%"),
    @("G12", "5"),
    @("E13", "0.11%"),
    @("G13", "5"),
    @("D14", "0.001581"),
    @("E14", "-0.70%"),
    @("G14", "5"),
    @("E15", "0.92%"),
    @("G15", "5"),
    @("D16", "0.0006464"),
    @("E16", "-1.03%"),
    @("G16", "5"),
    @("D17", "0.006256"),
    @("E17", "4.46%"),
    @("G17", "5"),
    @("D18", "3.459"),
    @("E18", "-0.36%"),
    @("G18", "5"),
    @("D19", "3.409"),
    @("E19", "0.41%"),
    @("G19", "5"),
    @("D20", "2.219"),
    @("E20", "-1.01%"),
    @("G20", "5"),
    @("D21", "0.3213"),
    @("E21", "1.76%"),
    @("G21", "5"),
    @("E22", "-2.94%"),
    @("G22", "5"),
    @("D23", "4.076"),
    @("E23", "1.57%"),
    @("G23", "5"),
    @("D24", "0.1567"),
    @("E24", "3.13%"),
    @("G24", "5"),
    @("D25", "0.001194"),
    @("E25", "0.69%"),
    @("G25", "5"),
    @("D26", "0.004132"),
    @("E26", "-4.73%"),
    @("G26", "5"),
    @("D27", "0.0001250"),
    @("E27", "5.95%"),
    @("G27", "5"),
    @("D28", "0.0001618"),
    @("E28", "-1.03%"),
    @("G28", "5"),
    @("G29", "5"),
    @("G30", "5"),
    @("G31", "5"),
    @("G32", "5"),
    @("G33", "5"),
    @("G34", "5"),
    @("G35", "5"),
    @("G36", "5"),
    @("G37", "5"),
    @("G38", "5"),
    @("G39", "5"),
    @("D40", "0.04219"),
    @("E40", "1.38%"),
    @("G40", "5"),
    @("D41", "0.006736"),
    @("E41", "0.43%"),
    @("G41", "5"),
    @("D42", "0.1241"),
    @("E42", "-12.12%"),
    @("G42", "5"),
    @("D43", "0.001980"),
    @("E43", "-8.31%"),
    @("G43", "5"),
    @("D44", "0.01211"),
    @("E44", "-2.44%"),
    @("G44", "5"),
    @("D45", "0.00005656"),
    @("E45", "1.94%"),
    @("G45", "5"),
    @("D46", "1.972"),
    @("E46", "26.20%"),
    @("G46", "5"),
    @("D47", "0.01307"),
    @("E47", "-29.33%"),
    @("G47", "5"),
    @("G48", "5"),
    @("G49", "5"),
    @("G50", "5"),
    @("G51", "5")
)

foreach ($u in $updates) {
    $cellRef = $u[0]
    $val = $u[1]
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}